$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first row, shifting the existing company list down by one.
$ws.Rows.Item(1).Insert() | Out-Null

# Add the new header/title text in the freshly inserted row.
$ws.Range("A1").Value = "Companies to screen"
$ws.Range("A1").Font.Bold = $true

# Reflect the author's last on-screen selection.
$ws.Range("B3").Select() | Out-Null

# Keep the print orientation explicit (portrait), as saved in the workbook.
$ws.PageSetup.Orientation = 1
